$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.800.70'
$ws.Range('E2').Value = '  +1.64%  '
$ws.Range('D3').Value = '3.115.69'
$ws.Range('E3').Value = '  +1.32%  '
$ws.Range('E4').Value = '  +0.05%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '532.90'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +2.47%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '138.37'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +1.70%  '
$ws.Range('E7').Value = '  -0.15%  '
$c = $ws.Range('D8')
$c.NumberFormat = "@"
$c.Value = '0.497'
$c.Style = "Normal"
$ws.Range('E8').Value = '  +10.31%  '
$ws.Range('E9').Value = '  +0.17%  '
$ws.Range('E10').Value = '  +1.51%  '
$ws.Range('E11').Value = '  +4.40%  '
$ws.Range('E12').Value = '  +3.53%  '
$ws.Range('D13').Value = '3.655.41'
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '25.66'
$c.Style = "Normal"
$ws.Range('E14').Value = '  +1.69%  '
$ws.Range('E15').Value = '  +3.63%  '
$ws.Range('D16').Value = '57.906.94'
$ws.Range('E16').Value = '  +1.63%  '
$ws.Range('D17').Value = '3.118.67'
$ws.Range('E17').Value = '  +1.52%  '
$c = $ws.Range('D18')
$c.NumberFormat = "@"
$c.Value = '6.14'
$c.Style = "Normal"
$ws.Range('E18').Value = '  +4.68%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.80'
$c.Style = "Normal"
$ws.Range('E19').Value = '  +2.95%  '
$c = $ws.Range('D20')
$c.NumberFormat = "@"
$c.Value = '8.11'
$c.Style = "Normal"
$ws.Range('E20').Value = '  +3.55%  '
$c = $ws.Range('D21')
$c.NumberFormat = "@"
$c.Value = '373.45'
$c.Style = "Normal"
$ws.Range('E21').Value = '  +7.89%  '
$ws.Range('E22').Value = '  +0.01%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '5.72'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -1.22%  '
$ws.Range('E24').Value = '  +1.61%  '
$c = $ws.Range('D25')
$c.NumberFormat = "@"
$c.Value = '0.507'
$c.Style = "Normal"
$ws.Range('E25').Value = '  +2.14%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -0.20%  '
$ws.Range('D28').Value = '0.0₃0883'
$ws.Range('E28').Value = '  +2.59%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '7.61'
$c.Style = "Normal"
$ws.Range('E29').Value = '  +4.85%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '6.14'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +4.72%  '
$ws.Range('E31').Value = '  +0.35%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '21.50'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +3.77%  '
$c = $ws.Range('D33')
$c.NumberFormat = "@"
$c.Value = '5.14'
$c.Style = "Normal"
$ws.Range('E33').Value = '  +4.68%  '
$ws.Range('E34').Value = '  +3.16%  '
$c = $ws.Range('D35')
$c.NumberFormat = "@"
$c.Value = '160.41'
$c.Style = "Normal"
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '6.18'
$c.Style = "Normal"
$ws.Range('E36').Value = '  +3.25%  '
$ws.Range('E37').Value = '  +6.30%  '
$ws.Range('E38').Value = '  -0.33%  '
$ws.Range('E39').Value = '  +4.31%  '
$c = $ws.Range('D40')
$c.NumberFormat = "@"
$c.Value = '0.0670'
$c.Style = "Normal"
$ws.Range('E40').Value = '  +2.97%  '
$ws.Range('D41').Value = '2.559.62'
$ws.Range('E41').Value = '  +7.63%  '
$ws.Range('E42').Value = '  +3.94%  '
$ws.Range('E43').Value = '  +5.05%  '
$ws.Range('E44').Value = '  +1.25%  '
$ws.Range('E45').Value = '  +2.84%  '
$ws.Range('E46').Value = '  +0.03%  '
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '0.980'
$c.Style = "Normal"
$ws.Range('E47').Value = '  +2.72%  '
$ws.Range('E48').Value = '  +3.53%  '
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '19.91'
$c.Style = "Normal"
$ws.Range('E49').Value = '  +1.51%  '
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '0.0950'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +6.88%  '
$ws.Range('E51').Value = '  -0.43%  '
